$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("L5").Value = 1.22
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 1.75
$ws.Range("O5").Value = 2.05
$ws.Range("H9").Value = 3.15
$ws.Range("I9").Value = 3.25
$ws.Range("L9").Value = 1.44
$ws.Range("M9").Value = 2.4
$ws.Range("N9").Value = 2.27
$ws.Range("O9").Value = 1.5
$ws.Range("P9").Value = 1.52
$ws.Range("Q9").Value = 2.22
$ws.Range("R9").Value = 2.02
$ws.Range("T9").Value = 6
$ws.Range("Z9").Value = 7.2
$ws.Range("AA9").Value = 6.2
$ws.Range("AB9").Value = 18.5
$ws.Range("AF9").Value = 15
$ws.Range("AG9").Value = 12
$ws.Range("AI9").Value = 35
$ws.Range("G12").Value = 3.6
$ws.Range("I12").Value = 2.62
$ws.Range("J12").Value = 1.19
$ws.Range("K12").Value = 4.15
$ws.Range("M12").Value = 2.02
$ws.Range("R12").Value = 2.25
$ws.Range("T12").Value = 6.8
$ws.Range("U12").Value = 17.5
$ws.Range("Z12").Value = 4.15
$ws.Range("AB12").Value = 18.5
$ws.Range("AH12").Value = 32
$ws.Range("G15").Value = 1.7
$ws.Range("H15").Value = 3.55
$ws.Range("I15").Value = 4.25
$ws.Range("N15").Value = 1.9
$ws.Range("O15").Value = 1.72
$ws.Range("T15").Value = 5.5
$ws.Range("U15").Value = 6.4
$ws.Range("V15").Value = 7.1
$ws.Range("W15").Value = 10.5
$ws.Range("X15").Value = 11.5
$ws.Range("Y15").Value = 23
$ws.Range("Z15").Value = 9.25
$ws.Range("AA15").Value = 6.1
$ws.Range("AB15").Value = 14
$ws.Range("AC15").Value = 65
$ws.Range("AD15").Value = 500
$ws.Range("AE15").Value = 9.5
$ws.Range("AF15").Value = 18.5
$ws.Range("AG15").Value = 12
$ws.Range("AH15").Value = 50
$ws.Range("AI15").Value = 32
$ws.Range("AJ15").Value = 40
$ws.Range("G16").Value = 3.15
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 2.18
$ws.Range("T16").Value = 7.6
$ws.Range("U16").Value = 13.5
$ws.Range("V16").Value = 9.25
$ws.Range("W16").Value = 32
$ws.Range("X16").Value = 23
$ws.Range("Y16").Value = 28
$ws.Range("AA16").Value = 5.2
$ws.Range("AB16").Value = 11.75
$ws.Range("AE16").Value = 5.8
$ws.Range("AF16").Value = 8.25
$ws.Range("AG16").Value = 7.6
$ws.Range("AH16").Value = 17
$ws.Range("AI16").Value = 15.5
$ws.Range("AJ16").Value = 25
$ws.Range("I18").Value = 3.85
$ws.Range("P18").Value = 1.39
$ws.Range("Q18").Value = 2.45
$ws.Range("T18").Value = 5.5
$ws.Range("U18").Value = 6.7
$ws.Range("W18").Value = 11.5
$ws.Range("AE18").Value = 8.5
$ws.Range("AG18").Value = 11.25
$ws.Range("AI18").Value = 30
$ws.Range("H22").Value = 3.55
$ws.Range("I22").Value = 3.2
$ws.Range("M22").Value = 4.5
$ws.Range("U22").Value = 13.5
$ws.Range("W22").Value = 22
$ws.Range("AA22").Value = 7.4
$ws.Range("AF22").Value = 21
$ws.Range("AG22").Value = 11.25
$ws.Range("P25").Value = 1.4
$ws.Range("T26").Value = 5.7
$ws.Range("Z26").Value = 11.75
$ws.Range("AB26").Value = 19.5
$ws.Range("AJ26").Value = 70
$ws.Range("N27").Value = 1.67
$ws.Range("L28").Value = 1.2
$ws.Range("M28").Value = 4.33
$ws.Range("N28").Value = 1.67
$ws.Range("O28").Value = 2.15
$ws.Range("P28").Value = 1.3
$ws.Range("Q28").Value = 3.4
$ws.Range("AB28").Value = 19
$ws.Range("AD28").Value = 251
$ws.Range("AH28").Value = 51
$ws.Range("P29").Value = 1.36
$ws.Range("G32").Value = 1.9
$ws.Range("H32").Value = 3.65
$ws.Range("I32").Value = 3.55
$ws.Range("M32").Value = 3.4
$ws.Range("N32").Value = 1.82
$ws.Range("O32").Value = 1.9
$ws.Range("P32").Value = 1.39
$ws.Range("Q32").Value = 2.77
$ws.Range("R32").Value = 1.75
$ws.Range("S32").Value = 1.98
$ws.Range("T32").Value = 7.7
$ws.Range("U32").Value = 9.25
$ws.Range("W32").Value = 16
$ws.Range("X32").Value = 14.5
$ws.Range("Y32").Value = 26
$ws.Range("AA32").Value = 7
$ws.Range("AB32").Value = 14.5
$ws.Range("AC32").Value = 65
$ws.Range("AD32").Value = 500
$ws.Range("AF32").Value = 19
$ws.Range("AG32").Value = 12
$ws.Range("AH32").Value = 50
$ws.Range("AI32").Value = 30
$ws.Range("AJ32").Value = 37
